$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture the text currently stored in the cells we still need,
#        before we start deleting / moving rows around. ---
$v_A10 = $ws.Range("A10").Value2   # Objetivos:
$v_A11 = $ws.Range("A11").Value2   # Objectives:
$v_A17 = $ws.Range("A17").Value2   # Programa resumido:
$v_B13 = $ws.Range("B13").Value2   # 471420 - Carlos Antonio Reis Pereira Baptista
$v_B14 = $ws.Range("B14").Value2   # 3480026 - Joao Paulo Pascon
$v_A18 = $ws.Range("A18").Value2   # Short syllabus:
$v_A19 = $ws.Range("A19").Value2   # Programa:
$v_B15 = $ws.Range("B15").Value2   # 5840793 - Sergio Schneider
$v_A20 = $ws.Range("A20").Value2   # Syllabus:
$v_A21 = $ws.Range("A21").Value2   # Avaliacao:
$v_A22 = $ws.Range("A22").Value2   # Metodo:
$v_B16 = $ws.Range("B16").Value2   # 7797767 - Viktor Pastoukhov
$v_A23 = $ws.Range("A23").Value2   # Criterio:
$v_B22 = $ws.Range("B22").Value2   # Os alunos serao avaliados...
$v_A24 = $ws.Range("A24").Value2   # Norma de recuperacao:
$v_B23 = $ws.Range("B23").Value2   # Nota Final (NF) = ...
$v_A25 = $ws.Range("A25").Value2   # Bibliografia:
$v_B24 = $ws.Range("B24").Value2   # Para a recuperacao sera realizada...
$v_A26 = $ws.Range("A26").Value2   # Requisitos:
$v_B27 = $ws.Range("B27").Value2   # LOM3099 - Estatica (Requisito fraco)

# --- 2. Remove old rows 10-27 completely (content, formatting, heights). ---
$ws.Range("A10:C27").EntireRow.Delete() | Out-Null

# --- 3. Rebuild rows 10-22 with the new layout. Styles are copied from the
#        existing A/B/C column formats (rows 1-9, untouched) so the cellXfs
#        table in styles.xml is not altered. ---

function Set-Cell($addr, $styleSrc, $value) {
    $ws.Range($styleSrc).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($addr).Value2 = $value
}

# Row 10 (height 60)
$ws.Rows.Item(10).RowHeight = 60
Set-Cell "A10" "A3" $v_A10
Set-Cell "B10" "B3" $v_B13
Set-Cell "C10" "C3" $v_B13

# Row 11 (height 60)
$ws.Rows.Item(11).RowHeight = 60
Set-Cell "A11" "A3" $v_A11

# Row 12 (height 60)
$ws.Rows.Item(12).RowHeight = 60
Set-Cell "A12" "A3" $v_A17
Set-Cell "B12" "B3" $v_B14
Set-Cell "C12" "C3" $v_B14

# Row 13 (height 60)
$ws.Rows.Item(13).RowHeight = 60
Set-Cell "A13" "A3" $v_A18

# Row 14 (height 120)
$ws.Rows.Item(14).RowHeight = 120
Set-Cell "A14" "A3" $v_A19
Set-Cell "B14" "B3" $v_B15
Set-Cell "C14" "C3" $v_B15

# Row 15 (height 120)
$ws.Rows.Item(15).RowHeight = 120
Set-Cell "A15" "A3" $v_A20

# Row 16 (default height)
Set-Cell "A16" "A3" $v_A21

# Row 17 (height 60)
$ws.Rows.Item(17).RowHeight = 60
Set-Cell "A17" "A3" $v_A22
Set-Cell "B17" "B3" $v_B16
Set-Cell "C17" "C3" $v_B16

# Row 18 (height 60)
$ws.Rows.Item(18).RowHeight = 60
Set-Cell "A18" "A3" $v_A23
Set-Cell "B18" "B3" $v_B22
Set-Cell "C18" "C3" $v_B22

# Row 19 (height 60)
$ws.Rows.Item(19).RowHeight = 60
Set-Cell "A19" "A3" $v_A24
Set-Cell "B19" "B3" $v_B23
Set-Cell "C19" "C3" $v_B23

# Row 20 (height 120)
$ws.Rows.Item(20).RowHeight = 120
Set-Cell "A20" "A3" $v_A25
Set-Cell "B20" "B3" $v_B24
Set-Cell "C20" "C3" $v_B24

# Row 21 (default height)
Set-Cell "A21" "A3" $v_A26

# Row 22 (height 30)
$ws.Rows.Item(22).RowHeight = 30
Set-Cell "B22" "B3" $v_B27
Set-Cell "C22" "C3" $v_B27

$ws.Range("A1").Select() | Out-Null
